# "task and table maj" -- update the Accueil Admin block and add a new
# "Type Gestion" task block underneath it.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Update existing cells ---
$ws.Range("B29").Value = "(ok) accueilAdmin.html (1h30)"
$ws.Range("D31").Value = "> accueil avec nav bar : variété - parcelles – cueilleurs – catDépenses – salKilo"
$ws.Range("D39").Value = "> modele_accueil.Admin.php"

# row 39/40 change from the default 15pt height to 13.8pt (matches rows 1-29 height)
$ws.Rows.Item(39).RowHeight = 13.8
$ws.Rows.Item(40).RowHeight = 13.8

# --- Insert new rows of content below (rows 40-41 replace/extend the old row 40) ---
$ws.Range("D40").Value = "> folder : type Gestion ; page : insert ou list"
$ws.Range("D41").Value = "> rediriger vers les bonnes pages de crud "

# --- New "Type Gestion" task block (rows 43-53), mirrors the Accueil Admin block ---
$ws.Range("B43").Value = "Type Gestion"
$ws.Range("B43").Font.Bold = $true
$ws.Rows.Item(43).RowHeight = 13.8

$ws.Range("C44").Value = "Affichage "
$ws.Range("D45").Value = "> formulaire de chq entité"

$ws.Range("C47").Value = "Métier "
$ws.Range("D48").Value = "> getListeVariete by Parcelle "
$ws.Range("D49").Value = "> fonction js : displayInTable (tabObject) "

$ws.Range("C51").Value = "Intégration "
$ws.Range("D52").Value = "> traitement AJAX pour chaque insertion "
$ws.Range("D53").Value = "> traitement AJAX pour list entité"

# --- Update the active selection to match the new working area ---
$ws.Range("D49").Select() | Out-Null
